$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9 (Total Cancer Risk (per million))
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("D9").Value = 54
$wsMeans.Range("E9").Value = 43
$wsMeans.Range("F9").Value = 42
$wsMeans.Range("G9").Value = 36

# Row 10 (Total Respiratory (hazard quotient))
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.43
$wsMeans.Range("F10").Value = 0.47
$wsMeans.Range("G10").Value = 0.43

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9 (Total Cancer Risk (per million))
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 5.4
$wsSD.Range("D9").Value = 17
$wsSD.Range("E9").Value = 11
$wsSD.Range("F9").Value = 8.3
$wsSD.Range("G9").Value = 8.6

# Row 10 (Total Respiratory (hazard quotient))
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.065
$wsSD.Range("F10").Value = 0.072
$wsSD.Range("G10").Value = 0.068
